$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").NumberFormat = "General"
$ws.Range("F2").Value = 2

$ws.Range("H2").NumberFormat = "General"
$ws.Range("H2").Value = 1

$ws.Range("O2").NumberFormat = "General"
$ws.Range("O2").Value = 1

$ws.Range("U2").NumberFormat = "General"
$ws.Range("U2").Value = 1

$ws.Range("V2").NumberFormat = "General"
$ws.Range("V2").Value = 1

$ws.Range("W2").NumberFormat = "General"
$ws.Range("W2").Value = 1

[void]$ws.Range("G5").Select()
